# Financials update: insert a new "most recent fiscal year" column (D)
# into the INS sheet, ahead of the existing year columns (which shift
# right from D:K to E:L), and populate the new column with the new
# year's reported figures for Income Statement, Balance Sheet, and
# Cash Flow Statement sections.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns D:K one column to the right by inserting a
# fresh column at D. Excel carries values/number-formats along for the
# ride, but the *new* column D cells come back blank with a generic
# style, so we still need to copy number formats across from the
# column that's now immediately to the right (E).
$ws.Range("D1").EntireColumn.Insert()

# Re-apply formats (incl. number formats) onto the new column D from
# column E, restricted to the three data blocks on the sheet (title
# rows such as 5/6/37/79 have no D:L cells at all and must stay that
# way).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Income Statement (rows 7-35) ---
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 20100
$ws.Range("D9").Value = 8500
$ws.Range("D10").Value = 11600
$ws.Range("D12").Value = 3400
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 400
$ws.Range("D15").Value = "NA"
$ws.Range("D17").Value = "NA"
$ws.Range("D18").Value = "NA"
$ws.Range("D20").Value = "NA"
$ws.Range("D21").Value = "NA"
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 6200
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = "NA"
$ws.Range("D27").Value = "NA"
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = "NA"
$ws.Range("D33").Value = "NA"
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = "NA"

# --- Balance Sheet (rows 38-77) ---
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 18900
$ws.Range("D42").Value = 300
$ws.Range("D43").Value = 4300
$ws.Range("D44").Value = "NA"
$ws.Range("D45").Value = 1200
$ws.Range("D46").Value = 24800
$ws.Range("D47").Value = 2500
$ws.Range("D48").Value = 1500
$ws.Range("D49").Value = "NA"
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 29300
$ws.Range("D57").Value = 500
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 3000
$ws.Range("D60").Value = 3300
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 100
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = "NA"
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 10900
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 25900
$ws.Range("D77").Value = 0

# --- Cash Flow Statement (rows 80-102) ---
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = "NA"
$ws.Range("D83").Value = 600
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 6700
$ws.Range("D91").Value = -900
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -1900
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 4900
